# Updates cryptos list figures (prices / 1h volume deltas) and restores
# the correct ImmutableX / NEARProtocol row ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$updates = @{
    "D2"  = "62.663.97"
    "E2"  = "  -1.00%  "
    "D3"  = "3.443.14"
    "E3"  = "  -1.24%  "
    "D4"  = "1.00"
    "E4"  = "  -0.01%  "
    "D5"  = "578.61"
    "E5"  = "  -1.04%  "
    "D6"  = "147.68"
    "E6"  = "  -0.22%  "
    "E7"  = "  +0.06%  "
    "D8"  = "0.479"
    "E8"  = "  +0.05%  "
    "D9"  = "7.98"
    "E9"  = "  +3.59%  "
    "E10" = "  -2.48%  "
    "E11" = "  +2.29%  "
    "D12" = "4.033.60"
    "E12" = "  -1.21%  "
    "E13" = "  +2.54%  "
    "D14" = "28.22"
    "E14" = "  -5.47%  "
    "D15" = "3.436.39"
    "E15" = "  -1.74%  "
    "E16" = "  -1.28%  "
    "D17" = "62.744.50"
    "E17" = "  -0.90%  "
    "D18" = "6.40"
    "E18" = "  +1.23%  "
    "D19" = "14.59"
    "E19" = "  +1.32%  "
    "D20" = "9.07"
    "E20" = "  -3.16%  "
    "D21" = "387.55"
    "E21" = "  -0.86%  "
    "D22" = "0.561"
    "E22" = "  -0.84%  "
    "D23" = "75.04"
    "E23" = "  -0.22%  "
    "E24" = "  +0.12%  "
    "E25" = "  -2.44%  "
    "D26" = "3.588.76"
    "E27" = "  +0.00%  "
    "D28" = "7.61"
    "E28" = "  -2.92%  "
    "E29" = "  +0.14%  "
    "D30" = "7.98"
    "E30" = "  -3.91%  "
    "E31" = "  -2.13%  "
    "E32" = "  +0.00%  "
    "D33" = "1.34"
    "E33" = "  -8.56%  "
    "E34" = "  -2.61%  "

    # Row 35 / 36: ImmutableX and NEARProtocol swap places.
    "B35" = "NEARProtocol"
    "C35" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D35" = "5.32"
    "E35" = "  -0.86%  "
    "B36" = "ImmutableX"
    "C36" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D36" = "1.62"
    "E36" = "  +2.61%  "

    "D37" = "32.07"
    "E37" = "  -0.03%  "
    "E38" = "  -2.31%  "
    "D39" = "170.52"
    "E39" = "  -0.49%  "
    "D40" = "3.479.19"
    "E40" = "  -1.22%  "
    "D41" = "0.0775"
    "E41" = "  +0.61%  "
    "D42" = "0.788"
    "E42" = "  -2.67%  "
    "D43" = "42.46"
    "E43" = "  -0.03%  "
    "E44" = "  -2.20%  "
    "E45" = "  -3.68%  "
    "E46" = "  -3.10%  "
    "D47" = "2.561.54"
    "E47" = "  -2.24%  "
    "E48" = "  +1.69%  "
    "E49" = "  -1.34%  "
    "D50" = "22.63"
    "E50" = "  -4.57%  "
    "E51" = "  -0.02%  "
}

# Some of the new "Price" figures (single decimal point, e.g. "1.00" or
# "578.61") parse cleanly as numbers, so Excel would otherwise silently
# convert them away from the text values used throughout this sheet.
# Force those specific cells to a text number format before writing them
# so they stay text, exactly like their neighbours (e.g. "62.663.97",
# which uses "." as a thousands separator and is never auto-converted).
$forceTextCells = @(
    "D4","D5","D6","D8","D9","D14","D18","D19","D20","D21","D22","D23",
    "D28","D30","D33","D35","D36","D37","D39","D41","D42","D43","D50"
)
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
